$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.053.50"
$ws.Range("E2").Value = "  -0.28%  "
$ws.Range("D3").Value = "1.829.07"
$ws.Range("E3").Value = "  -0.20%  "
$ws.Range("D4").Value = "'0.9987"
$ws.Range("E4").Value = "  -0.10%  "
$ws.Range("D5").Value = "'240.86"
$ws.Range("E5").Value = "  -0.23%  "
$ws.Range("D6").Value = "'0.6204"
$ws.Range("E6").Value = "  -6.37%  "
$ws.Range("E7").Value = "  -0.01%  "
$ws.Range("D8").Value = "'44.49"
$ws.Range("E8").Value = "  +6.04%  "
$ws.Range("D9").Value = "'0.07355"
$ws.Range("E9").Value = "  -0.78%  "
$ws.Range("D10").Value = "'0.2920"
$ws.Range("E10").Value = "  -0.47%  "
$ws.Range("D11").Value = "'22.69"
$ws.Range("E11").Value = "  +0.09%  "
$ws.Range("D12").Value = "'0.07683"
$ws.Range("E12").Value = "  -0.73%  "
$ws.Range("D13").Value = "1.832.78"
$ws.Range("E13").Value = "  +0.21%  "
$ws.Range("D14").Value = "'4.964"
$ws.Range("E14").Value = "  -0.41%  "
$ws.Range("D15").Value = "'0.6619"
$ws.Range("E15").Value = "  -1.06%  "
$ws.Range("D16").Value = "'81.94"
$ws.Range("E16").Value = "  -1.12%  "
$ws.Range("D17").Value = "'0.000009041"
$ws.Range("E17").Value = "  +7.81%  "
$ws.Range("D18").Value = "'6.021"
$ws.Range("E18").Value = "  -1.12%  "
$ws.Range("D19").Value = "29.060.06"
$ws.Range("E19").Value = "  -0.22%  "
$ws.Range("D20").Value = "2.077.69"
$ws.Range("E20").Value = "  +0.69%  "
$ws.Range("D21").Value = "'225.26"
$ws.Range("E21").Value = "  -0.83%  "
$ws.Range("D22").Value = "'12.36"
$ws.Range("E22").Value = "  -0.79%  "
$ws.Range("E23").Value = "  -0.08%  "
$ws.Range("D24").Value = "'7.142"
$ws.Range("E24").Value = "  -0.31%  "
$ws.Range("E25").Value = "  -0.03%  "
$ws.Range("D26").Value = "'159.74"
$ws.Range("E26").Value = "  -0.09%  "
$ws.Range("D27").Value = "'8.420"
$ws.Range("E27").Value = "  -2.33%  "
$ws.Range("D28").Value = "'0.1355"
$ws.Range("E28").Value = "  -3.33%  "
$ws.Range("D29").Value = "'17.79"
$ws.Range("E29").Value = "  -0.92%  "
$ws.Range("D30").Value = "'1.497"
$ws.Range("E30").Value = "  -0.94%  "
$ws.Range("D31").Value = "'4.036"
$ws.Range("E31").Value = "  -0.14%  "
$ws.Range("D32").Value = "'4.049"
$ws.Range("E32").Value = "  -1.49%  "
$ws.Range("D33").Value = "'1.200"
$ws.Range("E33").Value = "  +0.51%  "
$ws.Range("D34").Value = "'0.05241"
$ws.Range("E34").Value = "  -2.09%  "
$ws.Range("D35").Value = "'1.842"
$ws.Range("E35").Value = "  -1.64%  "
$ws.Range("D36").Value = "'1.152"
$ws.Range("E36").Value = "  +1.49%  "
$ws.Range("D37").Value = "'0.7311"
$ws.Range("E37").Value = "  -3.31%  "
$ws.Range("E38").Value = "  -1.02%  "
$ws.Range("D39").Value = "1.290.42"
$ws.Range("E39").Value = "  +0.85%  "
$ws.Range("D40").Value = "'2.751"
$ws.Range("E40").Value = "  +0.73%  "
$ws.Range("D41").Value = "'0.01784"
$ws.Range("E41").Value = "  -0.89%  "
$ws.Range("D42").Value = "'6.287"
$ws.Range("E42").Value = "  +5.26%  "
$ws.Range("D43").Value = "'0.9008"
$ws.Range("E43").Value = "  -3.04%  "
$ws.Range("D44").Value = "'0.9993"
$ws.Range("E44").Value = "  -0.09%  "
$ws.Range("D45").Value = "'101.87"
$ws.Range("E45").Value = "  -0.94%  "
$ws.Range("D46").Value = "1.974.34"
$ws.Range("E46").Value = "  +0.53%  "
$ws.Range("B47").Value = "Aave"
$ws.Range("C47").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D47").Value = "'63.99"
$ws.Range("E47").Value = "  +1.16%  "
$ws.Range("B48").Value = "Mantle"
$ws.Range("C48").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D48").Value = "'0.5114"
$ws.Range("E48").Value = "  -0.90%  "
$ws.Range("E49").Value = "  -0.38%  "
$ws.Range("D50").Value = "'1.714"
$ws.Range("E50").Value = "  -3.15%  "
$ws.Range("D51").Value = "'0.3965"
$ws.Range("E51").Value = "  -1.69%  "
